# Efna1-Epha1.xlsx: update with new TPM-derived NATMI stats.
# - Sending-cluster label for rows 12-16 changes from "MuSCs" to "Inflammatory-Mac"
#   (and all of its ligand/receptor/edge statistics are recomputed).
# - A brand-new sending-cluster block for "MuSCs" is appended as rows 17-21.
# - All numeric columns (E:T) across rows 2-16 are refreshed with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efna1"
$ws.Range("C2").Value = "Epha1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 14.02618833333333
$ws.Range("H2").Value = 42.078565
$ws.Range("I2").Value = 0.806325281849088
$ws.Range("J2").Value = 0.8172785134657441
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.310473333333333
$ws.Range("N2").Value = 9.931419999999999
$ws.Range("O2").Value = 0.1683295705132556
$ws.Range("P2").Value = 0.1783110568845311
$ws.Range("Q2").Value = 46.4333224458111
$ws.Range("R2").Value = 417.8999020122999
$ws.Range("S2").Value = 0.1357283883876368
$ws.Range("T2").Value = 0.1457297955050953

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efna1"
$ws.Range("C3").Value = "Epha1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 14.02618833333333
$ws.Range("H3").Value = 42.078565
$ws.Range("I3").Value = 0.806325281849088
$ws.Range("J3").Value = 0.8172785134657441
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.966196333333333
$ws.Range("N3").Value = 11.898589
$ws.Range("O3").Value = 0.2016715007605908
$ws.Range("P3").Value = 0.2136300730433972
$ws.Range("Q3").Value = 55.63061673830943
$ws.Range("R3").Value = 500.6755506447849
$ws.Range("S3").Value = 0.162612829691712
$ws.Range("T3").Value = 0.174595268528486

# Row 4: ECs -> Inflammatory-Mac
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efna1"
$ws.Range("C4").Value = "Epha1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 14.02618833333333
$ws.Range("H4").Value = 42.078565
$ws.Range("I4").Value = 0.806325281849088
$ws.Range("J4").Value = 0.8172785134657441
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.013056
$ws.Range("N4").Value = 15.039168
$ws.Range("O4").Value = 0.2549017854764673
$ws.Range("P4").Value = 0.2700167690767302
$ws.Range("Q4").Value = 70.31406758154665
$ws.Range("R4").Value = 632.82660823392
$ws.Range("S4").Value = 0.2055337540181483
$ws.Range("T4").Value = 0.2206789036418531

# Row 5: ECs -> MuSCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.02618833333333
$ws.Range("H5").Value = 42.078565
$ws.Range("I5").Value = 0.806325281849088
$ws.Range("J5").Value = 0.8172785134657441
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.3026905
$ws.Range("N5").Value = 6.605381
$ws.Range("O5").Value = 0.1679338322424817
$ws.Range("P5").Value = 0.1185945682727144
$ws.Range("Q5").Value = 46.32415895971082
$ws.Range("R5").Value = 277.9449537582649
$ws.Range("S5").Value = 0.1354092946149165
$ws.Range("T5").Value = 0.09692479246303569

# Row 6: ECs -> Resolving-Mac
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Efna1"
$ws.Range("C6").Value = "Epha1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.02618833333333
$ws.Range("H6").Value = 42.078565
$ws.Range("I6").Value = 0.806325281849088
$ws.Range("J6").Value = 0.8172785134657441
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.074201666666667
$ws.Range("N6").Value = 12.222605
$ws.Range("O6").Value = 0.2071633110072045
$ws.Range("P6").Value = 0.2194475327226272
$ws.Range("Q6").Value = 57.14551988464722
$ws.Range("R6").Value = 514.3096789618251
$ws.Range("S6").Value = 0.1670410151366744
$ws.Range("T6").Value = 0.179349753327274

# Row 7: FAPs -> ECs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efna1"
$ws.Range("C7").Value = "Epha1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.483247333333333
$ws.Range("H7").Value = 7.449742
$ws.Range("I7").Value = 0.1427547569137158
$ws.Range("J7").Value = 0.1446939568272663
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.310473333333333
$ws.Range("N7").Value = 9.931419999999999
$ws.Range("O7").Value = 0.1683295705132556
$ws.Range("P7").Value = 0.1783110568845311
$ws.Range("Q7").Value = 8.220724077071111
$ws.Range("R7").Value = 73.98651669363998
$ws.Range("S7").Value = 0.02402984692000998
$ws.Range("T7").Value = 0.02580053236667457

# Row 8: FAPs -> FAPs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Efna1"
$ws.Range("C8").Value = "Epha1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.483247333333333
$ws.Range("H8").Value = 7.449742
$ws.Range("I8").Value = 0.1427547569137158
$ws.Range("J8").Value = 0.1446939568272663
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.966196333333333
$ws.Range("N8").Value = 11.898589
$ws.Range("O8").Value = 0.2016715007605908
$ws.Range("P8").Value = 0.2136300730433972
$ws.Range("Q8").Value = 9.849046468226444
$ws.Range("R8").Value = 88.641418214038
$ws.Range("S8").Value = 0.02878956606750239
$ws.Range("T8").Value = 0.03091098056594707

# Row 9: FAPs -> Inflammatory-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Efna1"
$ws.Range("C9").Value = "Epha1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.483247333333333
$ws.Range("H9").Value = 7.449742
$ws.Range("I9").Value = 0.1427547569137158
$ws.Range("J9").Value = 0.1446939568272663
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.013056
$ws.Range("N9").Value = 15.039168
$ws.Range("O9").Value = 0.2549017854764673
$ws.Range("P9").Value = 0.2700167690767302
$ws.Range("Q9").Value = 12.44865794385067
$ws.Range("R9").Value = 112.037921494656
$ws.Range("S9").Value = 0.03638844242256522
$ws.Range("T9").Value = 0.03906979472742633

# Row 10: FAPs -> MuSCs
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Efna1"
$ws.Range("C10").Value = "Epha1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.483247333333333
$ws.Range("H10").Value = 7.449742
$ws.Range("I10").Value = 0.1427547569137158
$ws.Range("J10").Value = 0.1446939568272663
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.3026905
$ws.Range("N10").Value = 6.605381
$ws.Range("O10").Value = 0.1679338322424817
$ws.Range("P10").Value = 0.1185945682727144
$ws.Range("Q10").Value = 8.201397376950332
$ws.Range("R10").Value = 49.20838426170199
$ws.Range("S10").Value = 0.0239733533993642
$ws.Range("T10").Value = 0.01715991734160042

# Row 11: FAPs -> Resolving-Mac
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Efna1"
$ws.Range("C11").Value = "Epha1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 2.483247333333333
$ws.Range("H11").Value = 7.449742
$ws.Range("I11").Value = 0.1427547569137158
$ws.Range("J11").Value = 0.1446939568272663
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.074201666666667
$ws.Range("N11").Value = 12.222605
$ws.Range("O11").Value = 0.2071633110072045
$ws.Range("P11").Value = 0.2194475327226272
$ws.Range("Q11").Value = 10.11725042421222
$ws.Range("R11").Value = 91.05525381791
$ws.Range("S11").Value = 0.02957354810427398
$ws.Range("T11").Value = 0.03175273182561794

# Row 12: Inflammatory-Mac -> ECs
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("B12").Value = "Efna1"
$ws.Range("C12").Value = "Epha1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1863673333333333
$ws.Range("H12").Value = 0.559102
$ws.Range("I12").Value = 0.01071372271683668
$ws.Range("J12").Value = 0.01085925937435662
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.310473333333333
$ws.Range("N12").Value = 9.931419999999999
$ws.Range("O12").Value = 0.1683295705132556
$ws.Range("P12").Value = 0.1783110568845311
$ws.Range("Q12").Value = 0.6169640872044444
$ws.Range("R12").Value = 5.552676784839999
$ws.Range("S12").Value = 0.001803436343523228
$ws.Range("T12").Value = 0.001936326016024781

# Row 13: Inflammatory-Mac -> FAPs
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("B13").Value = "Efna1"
$ws.Range("C13").Value = "Epha1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1863673333333333
$ws.Range("H13").Value = 0.559102
$ws.Range("I13").Value = 0.01071372271683668
$ws.Range("J13").Value = 0.01085925937435662
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.966196333333333
$ws.Range("N13").Value = 11.898589
$ws.Range("O13").Value = 0.2016715007605908
$ws.Range("P13").Value = 0.2136300730433972
$ws.Range("Q13").Value = 0.7391694341197776
$ws.Range("R13").Value = 6.652524907078
$ws.Range("S13").Value = 0.002160652539037288
$ws.Range("T13").Value = 0.002319864373341001

# Row 14: Inflammatory-Mac -> Inflammatory-Mac
$ws.Range("A14").Value = "Inflammatory-Mac"
$ws.Range("B14").Value = "Efna1"
$ws.Range("C14").Value = "Epha1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1863673333333333
$ws.Range("H14").Value = 0.559102
$ws.Range("I14").Value = 0.01071372271683668
$ws.Range("J14").Value = 0.01085925937435662
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.013056
$ws.Range("N14").Value = 15.039168
$ws.Range("O14").Value = 0.2549017854764673
$ws.Range("P14").Value = 0.2700167690767302
$ws.Range("Q14").Value = 0.9342698785706666
$ws.Range("R14").Value = 8.408428907135999
$ws.Range("S14").Value = 0.002730947049621458
$ws.Range("T14").Value = 0.00293218213082997

# Row 15: Inflammatory-Mac -> MuSCs
$ws.Range("A15").Value = "Inflammatory-Mac"
$ws.Range("B15").Value = "Efna1"
$ws.Range("C15").Value = "Epha1"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1863673333333333
$ws.Range("H15").Value = 0.559102
$ws.Range("I15").Value = 0.01071372271683668
$ws.Range("J15").Value = 0.01085925937435662
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.3026905
$ws.Range("N15").Value = 6.605381
$ws.Range("O15").Value = 0.1679338322424817
$ws.Range("P15").Value = 0.1185945682727144
$ws.Range("Q15").Value = 0.6155136213103333
$ws.Range("R15").Value = 3.693081727862
$ws.Range("S15").Value = 0.001799196513421717
$ws.Range("T15").Value = 0.00128784917726325

# Row 16: Inflammatory-Mac -> Resolving-Mac
$ws.Range("A16").Value = "Inflammatory-Mac"
$ws.Range("B16").Value = "Efna1"
$ws.Range("C16").Value = "Epha1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1863673333333333
$ws.Range("H16").Value = 0.559102
$ws.Range("I16").Value = 0.01071372271683668
$ws.Range("J16").Value = 0.01085925937435662
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 4.074201666666667
$ws.Range("N16").Value = 12.222605
$ws.Range("O16").Value = 0.2071633110072045
$ws.Range("P16").Value = 0.2194475327226272
$ws.Range("Q16").Value = 0.759298100078889
$ws.Range("R16").Value = 6.83368290071
$ws.Range("S16").Value = 0.002219490271232989
$ws.Range("T16").Value = 0.002383037676897622

# Row 17: MuSCs -> ECs
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Efna1"
$ws.Range("C17").Value = "Epha1"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6993955000000001
$ws.Range("H17").Value = 1.398791
$ws.Range("I17").Value = 0.04020623852035952
$ws.Range("J17").Value = 0.02716827033263282
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 3.310473333333333
$ws.Range("N17").Value = 9.931419999999999
$ws.Range("O17").Value = 0.1683295705132556
$ws.Range("P17").Value = 0.1783110568845311
$ws.Range("Q17").Value = 2.315330152203333
$ws.Range("R17").Value = 13.89198091322
$ws.Range("S17").Value = 0.006767898862085631
$ws.Range("T17").Value = 0.004844402996736409

# Row 18: MuSCs -> FAPs
$ws.Range("A18").Value = "MuSCs"
$ws.Range("B18").Value = "Efna1"
$ws.Range("C18").Value = "Epha1"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.6993955000000001
$ws.Range("H18").Value = 1.398791
$ws.Range("I18").Value = 0.04020623852035952
$ws.Range("J18").Value = 0.02716827033263282
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 3.966196333333333
$ws.Range("N18").Value = 11.898589
$ws.Range("O18").Value = 0.2016715007605908
$ws.Range("P18").Value = 0.2136300730433972
$ws.Range("Q18").Value = 2.773939867649833
$ws.Range("R18").Value = 16.643639205899
$ws.Range("S18").Value = 0.008108452462339181
$ws.Range("T18").Value = 0.005803959575623111

# Row 19: MuSCs -> Inflammatory-Mac
$ws.Range("A19").Value = "MuSCs"
$ws.Range("B19").Value = "Efna1"
$ws.Range("C19").Value = "Epha1"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.6993955000000001
$ws.Range("H19").Value = 1.398791
$ws.Range("I19").Value = 0.04020623852035952
$ws.Range("J19").Value = 0.02716827033263282
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 5.013056
$ws.Range("N19").Value = 15.039168
$ws.Range("O19").Value = 0.2549017854764673
$ws.Range("P19").Value = 0.2700167690767302
$ws.Range("Q19").Value = 3.506108807648
$ws.Range("R19").Value = 21.036652845888
$ws.Range("S19").Value = 0.01024864198613236
$ws.Range("T19").Value = 0.007335888576620696

# Row 20: MuSCs -> MuSCs
$ws.Range("A20").Value = "MuSCs"
$ws.Range("B20").Value = "Efna1"
$ws.Range("C20").Value = "Epha1"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.6993955000000001
$ws.Range("H20").Value = 1.398791
$ws.Range("I20").Value = 0.04020623852035952
$ws.Range("J20").Value = 0.02716827033263282
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 3.3026905
$ws.Range("N20").Value = 6.605381
$ws.Range("O20").Value = 0.1679338322424817
$ws.Range("P20").Value = 0.1185945682727144
$ws.Range("Q20").Value = 2.30988687359275
$ws.Range("R20").Value = 9.239547494370999
$ws.Range("S20").Value = 0.006751987714779262
$ws.Range("T20").Value = 0.003222009290814983

# Row 21: MuSCs -> Resolving-Mac
$ws.Range("A21").Value = "MuSCs"
$ws.Range("B21").Value = "Efna1"
$ws.Range("C21").Value = "Epha1"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.6993955000000001
$ws.Range("H21").Value = 1.398791
$ws.Range("I21").Value = 0.04020623852035952
$ws.Range("J21").Value = 0.02716827033263282
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 4.074201666666667
$ws.Range("N21").Value = 12.222605
$ws.Range("O21").Value = 0.2071633110072045
$ws.Range("P21").Value = 0.2194475327226272
$ws.Range("Q21").Value = 2.849478311759167
$ws.Range("R21").Value = 17.096869870555
$ws.Range("S21").Value = 0.008329257495023083
$ws.Range("T21").Value = 0.005962009892837624

